# "show nested error messages" — move the "on error popup: display all
# nested errors" task (Id 48) from the Active sheet's Todo list to the
# Inactive sheet's Done list.

$wb = $excel.ActiveWorkbook

$wsActive = $wb.Worksheets.Item("Active")
$wsInactive = $wb.Worksheets.Item("Inactive")

# Remove the completed task from the Active sheet (row 4: Id 48).
$wsActive.Rows.Item(4).Delete()

# Insert the task as a new first data row on the Inactive sheet (row 2),
# matching the formatting of the existing data rows rather than the bold
# header row that Excel would otherwise copy down into the new row.
$wsInactive.Rows.Item(2).Insert()
$wsInactive.Rows.Item(2).Style = $wsInactive.Rows.Item(3).Style

$wsInactive.Range("A2").Value = 48
$wsInactive.Range("B2").Value = "on error popup:`ndisplay all nested errors"
$wsInactive.Range("C2").Value = "Done"
$wsInactive.Range("D2").Value = "Task"

$wsInactive.Range("E2").NumberFormat = "@"
$wsInactive.Range("E2").Value = "8/22/2018"
$wsInactive.Range("E2").Style = "Normal"

$wsInactive.Range("F2").NumberFormat = "@"
$wsInactive.Range("F2").Value = "8/22/2018"
$wsInactive.Range("F2").Style = "Normal"

# Undo the auto row-height growth Excel applies because of the embedded
# line break in the title, so the row reverts to the sheet's default.
$wsInactive.Rows.Item(2).AutoFit()
